# point source functions added
#
# - rename Sheet1 -> "Data"
# - convert the yr (B) column from a plain 4-digit year to a real date
#   (Jan 1 of that year), formatted m/d/yyyy, header "yr" -> "DATE"
# - add a new "Stations" worksheet (after "Data") with name/Lat/Long for
#   the two station ids (id1, id2) that already appear in the Data sheet
# - make "Stations" the active sheet/tab

$wb = $excel.ActiveWorkbook

# tighten the iterative-calculation convergence threshold
$excel.Iteration = $true
$excel.MaxIterations = 100
$excel.MaxChange = 0.0001

# ---------------------------------------------------------------------
# Data sheet (was "Sheet1")
# ---------------------------------------------------------------------
$data = $wb.Worksheets.Item(1)
$data.Name = "Data"

# header: yr -> DATE
$data.Range("B1").Value = "DATE"

# body: year number -> date serial (Jan 1 of that year), m/d/yyyy format
for ($r = 2; $r -le 41; $r++) {
    $yr = [int]$data.Cells.Item($r, 2).Value2
    $dt = (Get-Date -Year $yr -Month 1 -Day 1 -Hour 0 -Minute 0 -Second 0).Date
    $data.Cells.Item($r, 2).Value = $dt
}
$data.Range("B2:B41").NumberFormat = "m/d/yyyy"

$data.Activate() | Out-Null
$data.Range("I20").Select() | Out-Null

# ---------------------------------------------------------------------
# Stations sheet (new)
# ---------------------------------------------------------------------
$stations = $wb.Worksheets.Add($null, $data)
$stations.Name = "Stations"

$stations.Range("A1").Value = "name"
$stations.Range("B1").Value = "Lat"
$stations.Range("C1").Value = "Long"

$stations.Range("A2").Value = "id1"
$stations.Range("B2").Value = 52.62527
$stations.Range("C2").Value = 18.73177

$stations.Range("A3").Value = "id2"
$stations.Range("B3").Value = 52.57165
$stations.Range("C3").Value = 18.62141

$stations.Activate() | Out-Null
$stations.Range("D1").Select() | Out-Null
